# BOT; UPDATE DATA (#496)
# Applies the data refresh + formula/view changes described in the commit diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Workbook-level calculation option: switch reference style to R1C1
# (maps to <calcPr .../ refMode="R1C1"/> in the saved workbook.xml)
# ---------------------------------------------------------------------------
$excel.ReferenceStyle = -4150   # xlR1C1

# ---------------------------------------------------------------------------
# Sheet "all": add the SUM formula to column C (rows 2-17 as one shared
# formula, row 18 as its own formula) and refresh the daily counts that
# changed for 2020-12-01 .. 2020-12-06 (rows 12-18).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Activate()

$wsAll.Range("C2:C17").Formula = "=SUM(D2,G2:H2)"
$wsAll.Range("C18").Formula = "=SUM(D18,G18:H18)"

$wsAll.Range("D12").Value = 113
$wsAll.Range("E12").Value = 105

$wsAll.Range("D13").Value = 118
$wsAll.Range("E13").Value = 111

$wsAll.Range("D14").Value = 119
$wsAll.Range("E14").Value = 111

$wsAll.Range("D15").Value = 127
$wsAll.Range("E15").Value = 117

$wsAll.Range("D16").Value = 134
$wsAll.Range("E16").Value = 124

$wsAll.Range("D17").Value = 131
$wsAll.Range("E17").Value = 121

$wsAll.Range("D18").Value = 130
$wsAll.Range("E18").Value = 120
$wsAll.Range("H18").Value = 87

# cursor moved from A18 to B18
$wsAll.Range("B18").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "kobe": refresh the daily counts for rows 67-73 and scroll the
# frozen view so column G / row 49 is the new top-left of the scrollable
# pane, with the cursor resting on J73.
# ---------------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Activate()

$wsKobe.Range("F67").Value = 106
$wsKobe.Range("G67").Value = 99

$wsKobe.Range("F68").Value = 111
$wsKobe.Range("G68").Value = 105

$wsKobe.Range("F69").Value = 112
$wsKobe.Range("G69").Value = 105

$wsKobe.Range("F70").Value = 120
$wsKobe.Range("G70").Value = 111

$wsKobe.Range("F71").Value = 127
$wsKobe.Range("G71").Value = 118

$wsKobe.Range("F72").Value = 125
$wsKobe.Range("G72").Value = 116

$wsKobe.Range("B73").Value = 36
$wsKobe.Range("C73").Value = 1632
$wsKobe.Range("F73").Value = 124
$wsKobe.Range("G73").Value = 115
$wsKobe.Range("J73").Value = 82

$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 7
$wsKobe.Range("J73").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "other": scroll the frozen view so column B / row 33 is the new
# top-left of the scrollable pane, with the cursor resting on A48.
# ---------------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Activate()

$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 2
$wsOther.Range("A48").Select() | Out-Null

$wsAll.Activate() | Out-Null
